$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# The NALCO price table gets a brand-new "latest price" row inserted at
# the top (row 2); every existing data row shifts down by one, and a
# new last row (26) appears carrying what used to be the final row (25).
# We rewrite the full data block (rows 2-26) top to bottom with the
# post-edit values, then rebuild the Circular Link hyperlinks to match.
# ---------------------------------------------------------------------

# Extend formatting down to the newly-used row 26 by cloning row 25s
# look (borders / number formats / alignment) before writing any values.
$ws.Range("A25:F25").Copy($ws.Range("A26:F26"))

# Helper: write a string into a cell without letting Excels "looks like
# a date" auto-detection reinterpret it (the Circular Date column stores
# plain DD-MM-YYYY text, not real dates). We stage the literal text in a
# scratch cell pre-formatted as Text, then paste-special just the value
# into the destination so the destinations own number format/style is
# left completely untouched.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

function Set-TextValue($rng, [string]$text) {
    $scratch.Value2 = $text
    $scratch.Copy()
    $rng.PasteSpecial(-4163) | Out-Null  # xlPasteValues
}

$ws.Range("A2").Value2 = 25
$ws.Range("B2").Value2 = "ALUMINIUM INGOT"
$ws.Range("C2").Value2 = "IE07"
$ws.Range("D2").Value2 = 320.05
Set-TextValue $ws.Range("E2") "07-01-2026"
$ws.Range("F2").Value2 = "https://nalcoindia.com/wp-content/uploads/2026/01/Ingot-07-01-2026.pdf"

$ws.Range("A3").Value2 = 24
$ws.Range("B3").Value2 = "ALUMINIUM INGOT"
$ws.Range("C3").Value2 = "IE07"
$ws.Range("D3").Value2 = 307.25
Set-TextValue $ws.Range("E3") "01-01-2026"
$ws.Range("F3").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-01-2026.pdf"

$ws.Range("A4").Value2 = 23
$ws.Range("B4").Value2 = "ALUMINIUM INGOT"
$ws.Range("C4").Value2 = "IE07"
$ws.Range("D4").Value2 = 301.65
Set-TextValue $ws.Range("E4") "24-12-2025"
$ws.Range("F4").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-24-12-2025.pdf"

$ws.Range("A5").Value2 = 22
$ws.Range("B5").Value2 = "ALUMINIUM INGOT"
$ws.Range("C5").Value2 = "IE07"
$ws.Range("D5").Value2 = 296.05
Set-TextValue $ws.Range("E5") "05-12-2025"
$ws.Range("F5").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-05-12-2025.pdf"

$ws.Range("A6").Value2 = 20
$ws.Range("B6").Value2 = "ALUMINIUM INGOT"
$ws.Range("C6").Value2 = "IE07"
$ws.Range("D6").Value2 = 290.95
Set-TextValue $ws.Range("E6") "27-11-2025"
$ws.Range("F6").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-27-11-2025.pdf"

$ws.Range("A7").Value2 = 19
$ws.Range("B7").Value2 = "ALUMINIUM INGOT"
$ws.Range("C7").Value2 = "IE07"
$ws.Range("D7").Value2 = 283.55
Set-TextValue $ws.Range("E7") "22-11-2025"
$ws.Range("F7").Value2 = "https://nalcoindia.com/wp-content/uploads/2025/11/Ingot-22-11-2025.pdf"

$ws.Range("A8").Value2 = 18
$ws.Range("B8").Value2 = "ALUMINIUM INGOT"
$ws.Range("C8").Value2 = "IE07"
$ws.Range("D8").Value2 = 281.95
Set-TextValue $ws.Range("E8") "19-11-2025"
$ws.Range("F8").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-19-11-2025.pdf"

$ws.Range("A9").Value2 = 17
$ws.Range("B9").Value2 = "ALUMINIUM INGOT"
$ws.Range("C9").Value2 = "IE07"
$ws.Range("D9").Value2 = 292.65
Set-TextValue $ws.Range("E9") "07-11-2025"
$ws.Range("F9").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-11-2025.pdf"

$ws.Range("A10").Value2 = 21
$ws.Range("B10").Value2 = "ALUMINIUM INGOT"
$ws.Range("C10").Value2 = "IE07"
$ws.Range("D10").Value2 = 296.05
Set-TextValue $ws.Range("E10") "02-11-2025"
$ws.Range("F10").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf"

$ws.Range("A11").Value2 = 16
$ws.Range("B11").Value2 = "ALUMINIUM INGOT"
$ws.Range("C11").Value2 = "IE07"
$ws.Range("D11").Value2 = 297.15
Set-TextValue $ws.Range("E11") "01-11-2025"
$ws.Range("F11").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf"

$ws.Range("A12").Value2 = 15
$ws.Range("B12").Value2 = "ALUMINIUM INGOT"
$ws.Range("C12").Value2 = "IE07"
$ws.Range("D12").Value2 = 294.05
Set-TextValue $ws.Range("E12") "30-10-2025"
$ws.Range("F12").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-30-10-2025.pdf"

$ws.Range("A13").Value2 = 14
$ws.Range("B13").Value2 = "ALUMINIUM INGOT"
$ws.Range("C13").Value2 = "IE07"
$ws.Range("D13").Value2 = 288.55
Set-TextValue $ws.Range("E13") "25-10-2025"
$ws.Range("F13").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf"

$ws.Range("A14").Value2 = 13
$ws.Range("B14").Value2 = "ALUMINIUM INGOT"
$ws.Range("C14").Value2 = "IE07"
$ws.Range("D14").Value2 = 282.45
Set-TextValue $ws.Range("E14") "17-10-2025"
$ws.Range("F14").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf"

$ws.Range("A15").Value2 = 12
$ws.Range("B15").Value2 = "ALUMINIUM INGOT"
$ws.Range("C15").Value2 = "IE07"
$ws.Range("D15").Value2 = 285.05
Set-TextValue $ws.Range("E15") "14-10-2025"
$ws.Range("F15").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf"

$ws.Range("A16").Value2 = 11
$ws.Range("B16").Value2 = "ALUMINIUM INGOT"
$ws.Range("C16").Value2 = "IE07"
$ws.Range("D16").Value2 = 282.85
Set-TextValue $ws.Range("E16") "09-10-2025"
$ws.Range("F16").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf"

$ws.Range("A17").Value2 = 10
$ws.Range("B17").Value2 = "ALUMINIUM INGOT"
$ws.Range("C17").Value2 = "IE07"
$ws.Range("D17").Value2 = 277.95
Set-TextValue $ws.Range("E17") "01-10-2025"
$ws.Range("F17").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf"

$ws.Range("A18").Value2 = 9
$ws.Range("B18").Value2 = "ALUMINIUM INGOT"
$ws.Range("C18").Value2 = "IE07"
$ws.Range("D18").Value2 = 274.95
Set-TextValue $ws.Range("E18") "30-09-2025"
$ws.Range("F18").Value2 = "https://nalcoindia.com/wp-content/uploads/2025/09/INGOT-30-09-2025.pdf"

$ws.Range("A19").Value2 = 8
$ws.Range("B19").Value2 = "ALUMINIUM INGOT"
$ws.Range("C19").Value2 = "IE07"
$ws.Range("D19").Value2 = 270.25
Set-TextValue $ws.Range("E19") "25-09-2025"
$ws.Range("F19").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf"

$ws.Range("A20").Value2 = 7
$ws.Range("B20").Value2 = "ALUMINIUM INGOT"
$ws.Range("C20").Value2 = "IE07"
$ws.Range("D20").Value2 = 275.25
Set-TextValue $ws.Range("E20") "20-09-2025"
$ws.Range("F20").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf"

$ws.Range("A21").Value2 = 6
$ws.Range("B21").Value2 = "ALUMINIUM INGOT"
$ws.Range("C21").Value2 = "IE07"
$ws.Range("D21").Value2 = 278.95
Set-TextValue $ws.Range("E21") "17-09-2025"
$ws.Range("F21").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf"

$ws.Range("A22").Value2 = 5
$ws.Range("B22").Value2 = "ALUMINIUM INGOT"
$ws.Range("C22").Value2 = "IE07"
$ws.Range("D22").Value2 = 272.05
Set-TextValue $ws.Range("E22") "01-09-2025"
$ws.Range("F22").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf"

$ws.Range("A23").Value2 = 4
$ws.Range("B23").Value2 = "ALUMINIUM INGOT"
$ws.Range("C23").Value2 = "IE07"
$ws.Range("D23").Value2 = 271.05
Set-TextValue $ws.Range("E23") "28-08-2025"
$ws.Range("F23").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf"

$ws.Range("A24").Value2 = 3
$ws.Range("B24").Value2 = "ALUMINIUM INGOT"
$ws.Range("C24").Value2 = "IE07"
$ws.Range("D24").Value2 = 264.35
Set-TextValue $ws.Range("E24") "21-08-2025"
$ws.Range("F24").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"

$ws.Range("A25").Value2 = 2
$ws.Range("B25").Value2 = "ALUMINIUM INGOT"
$ws.Range("C25").Value2 = "IE07"
$ws.Range("D25").Value2 = 269.45
Set-TextValue $ws.Range("E25") "15-08-2025"
$ws.Range("F25").Value2 = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"

$ws.Range("A26").Value2 = 1
$ws.Range("B26").Value2 = "ALUMINIUM INGOT"
$ws.Range("C26").Value2 = "IE07"
$ws.Range("D26").Value2 = 268.25
Set-TextValue $ws.Range("E26") "07-08-2025"
$ws.Range("F26").Value2 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"

# Clean up the scratch cell so it does not end up inside the used range.
$scratch.Clear() | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Rebuild the Circular Link hyperlinks (F2:F26) so the relationship
# targets line up with the shifted rows. Cell text was already written
# above, so Hyperlinks.Add below only attaches the link target/rel -
# it will not touch the already-correct display text. A per-range
# Delete() is not reliably scoped in this host, so clear everything
# via the first hyperlinks Range, then re-add in row order so the
# relationship ids come out sequential again.
# ---------------------------------------------------------------------
if ($ws.Hyperlinks.Count -gt 0) {
    $ws.Hyperlinks.Item(1).Range.Hyperlinks.Delete()
}

$ws.Hyperlinks.Add($ws.Range("F2"), "https://nalcoindia.com/wp-content/uploads/2026/01/Ingot-07-01-2026.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-01-2026.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-24-12-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-05-12-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-27-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://nalcoindia.com/wp-content/uploads/2025/11/Ingot-22-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-19-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-02-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-11-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-30-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-25-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-17-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-14-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F16"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-09-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F17"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-01-10-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F18"), "https://nalcoindia.com/wp-content/uploads/2025/09/INGOT-30-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F19"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-25-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F20"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-20-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F21"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-17-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F22"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-01-09-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F23"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F24"), "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F25"), "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F26"), "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf") | Out-Null

Write-Output ("Final dimension: " + $ws.UsedRange.Address())
Write-Output ("Final hyperlink count: " + $ws.Hyperlinks.Count)
